$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume(1h) (E) columns to remain text so that
# values such as "0.9110" or "2.22%" are stored verbatim instead of being
# auto-converted into numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "289.18"
$ws.Range("D3").Value = "31.06"
$ws.Range("E3").Value = "2.22%"
$ws.Range("D4").Value = "4.952"
$ws.Range("E4").Value = "-0.01%"
$ws.Range("D5").Value = "0.07358"
$ws.Range("E5").Value = "1.74%"
$ws.Range("D6").Value = "2.348"
$ws.Range("E6").Value = "30.45%"
$ws.Range("D7").Value = "7.733"
$ws.Range("E7").Value = "2.23%"
$ws.Range("D8").Value = "3.726"
$ws.Range("E8").Value = "0.21%"
$ws.Range("D9").Value = "0.9110"
$ws.Range("E9").Value = "1.15%"
$ws.Range("D10").Value = "0.09266"
$ws.Range("E10").Value = "18.95%"
$ws.Range("D11").Value = "0.1701"
$ws.Range("D12").Value = "0.08167"
$ws.Range("E12").Value = "3.55%"
$ws.Range("D13").Value = "0.03114"
$ws.Range("E13").Value = "2.52%"
$ws.Range("D14").Value = "0.09974"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").Value = "0.32%"
$ws.Range("D16").Value = "0.005719"
$ws.Range("E16").Value = "-1.15%"
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").Value = "-0.03%"
$ws.Range("D18").Value = "2.106"
$ws.Range("E18").Value = "1.67%"
$ws.Range("D19").Value = "0.3325"
$ws.Range("E19").Value = "0.64%"
$ws.Range("E20").Value = "-0.85%"
$ws.Range("D21").Value = "4.179"
$ws.Range("E21").Value = "5.40%"
$ws.Range("D23").Value = "0.04525"
$ws.Range("E23").Value = "0.62%"
$ws.Range("D24").Value = "0.001212"
$ws.Range("E24").Value = "-0.10%"
$ws.Range("D25").Value = "0.004180"
$ws.Range("E25").Value = "-9.76%"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "0.07%"
$ws.Range("E39").Value = "0.53%"
$ws.Range("D40").Value = "0.04468"
$ws.Range("E40").Value = "3.11%"
$ws.Range("D41").Value = "0.007397"
$ws.Range("E41").Value = "1.04%"
$ws.Range("D42").Value = "0.009876"
$ws.Range("E42").Value = "-1.85%"
$ws.Range("D43").Value = "0.1330"
$ws.Range("E43").Value = "2.07%"
$ws.Range("D44").Value = "0.002240"
$ws.Range("E44").Value = "9.35%"
$ws.Range("D45").Value = "0.008762"
$ws.Range("E45").Value = "-7.03%"
$ws.Range("D46").Value = "0.00006107"
$ws.Range("E46").Value = "3.42%"
$ws.Range("E47").Value = "0.00%"
$ws.Range("D48").Value = "2.606"
$ws.Range("E48").Value = "15.58%"
$ws.Range("E50").Value = "0.00%"
$ws.Range("E51").Value = "0.00%"
